{"js": "// Update each two-digit \u00f7 one-digit division answer cell in the table,\n// replacing the old 'quotient, remainder' text with the new value for\n// each cell, in document order (matches the authored diff 1:1).\nconst replacements = [\n  [\"90\u00f78=11, 2\", \"73\u00f74=18, 1\"],\n  [\"21\u00f78=2, 5\", \"69\u00f79=7, 6\"],\n  [\"25\u00f73=8, 1\", \"27\u00f79=3, 0\"],\n  [\"47\u00f79=5, 2\", \"88\u00f75=17, 3\"],\n  [\"53\u00f74=13, 1\", \"50\u00f76=8, 2\"],\n  [\"59\u00f76=9, 5\", \"76\u00f77=10, 6\"],\n  [\"93\u00f77=13, 2\", \"24\u00f78=3, 0\"],\n  [\"68\u00f78=8, 4\", \"69\u00f72=34, 1\"],\n  [\"22\u00f73=7, 1\", \"21\u00f76=3, 3\"],\n  [\"46\u00f72=23, 0\", \"98\u00f79=10, 8\"],\n  [\"20\u00f73=6, 2\", \"95\u00f77=13, 4\"],\n  [\"84\u00f78=10, 4\", \"55\u00f77=7, 6\"],\n  [\"90\u00f72=45, 0\", \"47\u00f78=5, 7\"],\n  [\"82\u00f77=11, 5\", \"44\u00f76=7, 2\"],\n  [\"48\u00f75=9, 3\", \"54\u00f73=18, 0\"],\n  [\"31\u00f79=3, 4\", \"76\u00f79=8, 4\"],\n  [\"91\u00f79=10, 1\", \"49\u00f73=16, 1\"],\n  [\"65\u00f79=7, 2\", \"34\u00f78=4, 2\"],\n  [\"21\u00f73=7, 0\", \"96\u00f79=10, 6\"],\n  [\"95\u00f72=47, 1\", \"37\u00f75=7, 2\"],\n  [\"59\u00f79=6, 5\", \"33\u00f74=8, 1\"],\n  [\"89\u00f76=14, 5\", \"25\u00f78=3, 1\"],\n  [\"60\u00f77=8, 4\", \"52\u00f75=10, 2\"],\n  [\"95\u00f75=19, 0\", \"68\u00f74=17, 0\"],\n  [\"22\u00f74=5, 2\", \"59\u00f77=8, 3\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  // Each old value is unique in the document, so only the first hit matters;\n  // insertText(..., replace) swaps the text in place and keeps run formatting.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Apply each two-digit-division answer cell text update, matching the diff\n# exactly (old text -> new text), in document order using Find & Replace.\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceOne = 1\n\n$replacements = @(\n    ,@('90\u00f78=11, 2', '73\u00f74=18, 1')\n    ,@('21\u00f78=2, 5', '69\u00f79=7, 6')\n    ,@('25\u00f73=8, 1', '27\u00f79=3, 0')\n    ,@('47\u00f79=5, 2', '88\u00f75=17, 3')\n    ,@('53\u00f74=13, 1', '50\u00f76=8, 2')\n    ,@('59\u00f76=9, 5', '76\u00f77=10, 6')\n    ,@('93\u00f77=13, 2', '24\u00f78=3, 0')\n    ,@('68\u00f78=8, 4', '69\u00f72=34, 1')\n    ,@('22\u00f73=7, 1', '21\u00f76=3, 3')\n    ,@('46\u00f72=23, 0', '98\u00f79=10, 8')\n    ,@('20\u00f73=6, 2', '95\u00f77=13, 4')\n    ,@('84\u00f78=10, 4', '55\u00f77=7, 6')\n    ,@('90\u00f72=45, 0', '47\u00f78=5, 7')\n    ,@('82\u00f77=11, 5', '44\u00f76=7, 2')\n    ,@('48\u00f75=9, 3', '54\u00f73=18, 0')\n    ,@('31\u00f79=3, 4', '76\u00f79=8, 4')\n    ,@('91\u00f79=10, 1', '49\u00f73=16, 1')\n    ,@('65\u00f79=7, 2', '34\u00f78=4, 2')\n    ,@('21\u00f73=7, 0', '96\u00f79=10, 6')\n    ,@('95\u00f72=47, 1', '37\u00f75=7, 2')\n    ,@('59\u00f79=6, 5', '33\u00f74=8, 1')\n    ,@('89\u00f76=14, 5', '25\u00f78=3, 1')\n    ,@('60\u00f77=8, 4', '52\u00f75=10, 2')\n    ,@('95\u00f75=19, 0', '68\u00f74=17, 0')\n    ,@('22\u00f74=5, 2', '59\u00f77=8, 3')\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceOne) | Out-Null\n}\n"}
